$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# columns corresponding to each value in the per-row arrays below:
# E, G, H, I, J, K, M, N, O, P, Q, R, S, T
$cols = @(5,7,8,9,10,11,13,14,15,16,17,18,19,20)

$rowsData = @{}
$rowsData[2] = @(3; 95.32301533333333; 285.969046; 0.2086259240761469; 0.2086259240761469; 3; 13.38628533333333; 40.158856; 0.0790701185256053; 0.0790701185256053; 1276.021082085708; 11484.18973877138; 0.01649607654421487; 0.01649607654421487)
$rowsData[3] = @(3; 95.32301533333333; 285.969046; 0.2086259240761469; 0.2086259240761469; 3; 27.29177033333334; 81.87531100000001; 0.1612070459649248; 0.1612070459649248; 2601.533841958145; 23413.80457762331; 0.03363196893201833; 0.03363196893201833)
$rowsData[4] = @(3; 95.32301533333333; 285.969046; 0.2086259240761469; 0.2086259240761469; 3; 40.05126633333333; 120.153799; 0.2365748448668895; 0.2365748448668895; 3817.80747481175; 34360.26727330575; 0.04935564562352593; 0.04935564562352592)
$rowsData[5] = @(3; 95.32301533333333; 285.969046; 0.2086259240761469; 0.2086259240761469; 3; 40.10401433333334; 120.312043; 0.2368864167860689; 0.2368864167860689; 3822.835573224554; 34405.52015902098; 0.04942064760308091; 0.0494206476030809)
$rowsData[6] = @(3; 95.32301533333333; 285.969046; 0.2086259240761469; 0.2086259240761469; 3; 29.373837; 88.121511; 0.1735053986454552; 0.1735053986454552; 2800.002714749834; 25200.02443274851; 0.03619772412460835; 0.03619772412460835)
$rowsData[7] = @(3; 95.32301533333333; 285.969046; 0.2086259240761469; 0.2086259240761469; 3; 19.089213; 57.267639; 0.1127561752110562; 0.1127561752110562; 1819.641343500266; 16376.77209150239; 0.02352386124869854; 0.02352386124869853)
$rowsData[8] = @(3; 87.85366600000002; 263.560998; 0.1922783515464171; 0.1922783515464171; 3; 13.38628533333333; 40.158856; 0.0790701185256053; 0.0790701185256053; 1176.034240655366; 10584.30816589829; 0.01520347204668321; 0.0152034720466832)
$rowsData[9] = @(3; 87.85366600000002; 263.560998; 0.1922783515464171; 0.1922783515464171; 3; 27.29177033333334; 81.87531100000001; 0.1612070459649248; 0.1612070459649248; 2397.682075413376; 21579.13867872038; 0.03099662505580323; 0.03099662505580323)
$rowsData[10] = @(3; 87.85366600000002; 263.560998; 0.1922783515464171; 0.1922783515464171; 3; 40.05126633333333; 120.153799; 0.2365748448668895; 0.2365748448668895; 3518.650575325712; 31667.8551779314; 0.04548822118835488; 0.04548822118835487)
$rowsData[11] = @(3; 87.85366600000002; 263.560998; 0.1922783515464171; 0.1922783515464171; 3; 40.10401433333334; 120.312043; 0.2368864167860689; 0.2368864167860689; 3523.284680499881; 31709.56212449892; 0.04554812972336284; 0.04554812972336283)
$rowsData[12] = @(3; 87.85366600000002; 263.560998; 0.1922783515464171; 0.1922783515464171; 3; 29.373837; 88.121511; 0.1735053986454552; 0.1735053986454552; 2580.599264936443; 23225.39338442798; 0.03336133203595208; 0.03336133203595208)
$rowsData[13] = @(3; 87.85366600000002; 263.560998; 0.1922783515464171; 0.1922783515464171; 3; 19.089213; 57.267639; 0.1127561752110562; 0.1127561752110562; 1677.057343104858; 15093.51608794372; 0.02168057149626087; 0.02168057149626087)
$rowsData[14] = @(3; 89.00545999999999; 267.01638; 0.1947991917312883; 0.1947991917312882; 3; 13.38628533333333; 40.158856; 0.0790701185256053; 0.0790701185256053; 1191.452483784587; 10723.07235406128; 0.01540279517888508; 0.01540279517888507)
$rowsData[15] = @(3; 89.00545999999999; 267.01638; 0.1947991917312883; 0.1947991917312882; 3; 27.29177033333334; 81.87531100000001; 0.1612070459649248; 0.1612070459649248; 2429.116572732687; 21862.04915459418; 0.03140300225535599; 0.03140300225535598)
$rowsData[16] = @(3; 89.00545999999999; 267.01638; 0.1947991917312883; 0.1947991917312882; 3; 40.05126633333333; 120.153799; 0.2365748448668895; 0.2365748448668895; 3564.781383580846; 32083.03245222761; 0.04608458856402499; 0.04608458856402498)
$rowsData[17] = @(3; 89.00545999999999; 267.01638; 0.1947991917312883; 0.1947991917312882; 3; 40.10401433333334; 120.312043; 0.2368864167860689; 0.2368864167860689; 3569.476243584927; 32125.28619226434; 0.0461452825220473; 0.04614528252204728)
$rowsData[18] = @(3; 89.00545999999999; 267.01638; 0.1947991917312883; 0.1947991917312882; 3; 29.373837; 88.121511; 0.1735053986454552; 0.1735053986454552; 2614.431874150019; 23529.88686735018; 0.03379871141714964; 0.03379871141714963)
$rowsData[19] = @(3; 89.00545999999999; 267.01638; 0.1947991917312883; 0.1947991917312882; 3; 19.089213; 57.267639; 0.1127561752110562; 0.1127561752110562; 1699.04418410298; 15291.39765692682; 0.02196481179382527; 0.02196481179382527)
$rowsData[20] = @(3; 68.108298; 204.324894; 0.1490632305096083; 0.1490632305096083; 3; 13.38628533333333; 40.158856; 0.0790701185256053; 0.0790701185256053; 911.7171105956961; 8205.453995361266; 0.01178644730420435; 0.01178644730420435)
$rowsData[21] = @(3; 68.108298; 204.324894; 0.1490632305096083; 0.1490632305096083; 3; 27.29177033333334; 81.87531100000001; 0.1612070459649248; 0.1612070459649248; 1858.796026810226; 16729.16424129204; 0.0240300430524426; 0.0240300430524426)
$rowsData[22] = @(3; 68.108298; 204.324894; 0.1490632305096083; 0.1490632305096083; 3; 40.05126633333333; 120.153799; 0.2365748448668895; 0.2365748448668895; 2727.823582708034; 24550.41224437231; 0.03526461063316797; 0.03526461063316796)
$rowsData[23] = @(3; 68.108298; 204.324894; 0.1490632305096083; 0.1490632305096083; 3; 40.10401433333334; 120.312043; 0.2368864167860689; 0.2368864167860689; 2731.416159210939; 24582.74543289845; 0.03531105454997693; 0.03531105454997692)
$rowsData[24] = @(3; 68.108298; 204.324894; 0.1490632305096083; 0.1490632305096083; 3; 29.373837; 88.121511; 0.1735053986454552; 0.1735053986454552; 2000.602043799426; 18005.41839419484; 0.02586327523294897; 0.02586327523294896)
$rowsData[25] = @(3; 68.108298; 204.324894; 0.1490632305096083; 0.1490632305096083; 3; 19.089213; 57.267639; 0.1127561752110562; 0.1127561752110562; 1300.133807589474; 11701.20426830527; 0.01680779973686745; 0.01680779973686745)
$rowsData[26] = @(3; 72.54370233333333; 217.631107; 0.1587706482245977; 0.1587706482245977; 3; 13.38628533333333; 40.158856; 0.0790701185256053; 0.0790701185256053; 971.090698570399; 8739.816287133592; 0.01255401397350613; 0.01255401397350613)
$rowsData[27] = @(3; 72.54370233333333; 217.631107; 0.1587706482245977; 0.1587706482245977; 3; 27.29177033333334; 81.87531100000001; 0.1612070459649248; 0.1612070459649248; 1979.846063211031; 17818.61456889928; 0.02559494718622363; 0.02559494718622363)
$rowsData[28] = @(3; 72.54370233333333; 217.631107; 0.1587706482245977; 0.1587706482245977; 3; 40.05126633333333; 120.153799; 0.2365748448668895; 0.2365748448668895; 2905.467142958388; 26149.20428662549; 0.0375611414731497; 0.03756114147314969)
$rowsData[29] = @(3; 72.54370233333333; 217.631107; 0.1587706482245977; 0.1587706482245977; 3; 40.10401433333334; 120.312043; 0.2368864167860689; 0.2368864167860689; 2909.293678169067; 26183.6431035216; 0.03761060994872639; 0.03761060994872638)
$rowsData[30] = @(3; 72.54370233333333; 217.631107; 0.1587706482245977; 0.1587706482245977; 3; 29.373837; 88.121511; 0.1735053986454552; 0.1735053986454552; 2130.886887715853; 19177.98198944268; 0.02754756461340617; 0.02754756461340616)
$rowsData[31] = @(3; 72.54370233333333; 217.631107; 0.1587706482245977; 0.1587706482245977; 3; 19.089213; 57.267639; 0.1127561752110562; 0.1127561752110562; 1384.802185649597; 12463.21967084637; 0.01790237102958572; 0.01790237102958571)
$rowsData[32] = @(3; 44.074633; 132.223899; 0.09646265391194166; 0.09646265391194163; 3; 13.38628533333333; 40.158856; 0.0790701185256053; 0.0790701185256053; 589.9956132999494; 5309.960519699544; 0.007627313478111671; 0.007627313478111668)
$rowsData[33] = @(3; 44.074633; 132.223899; 0.09646265391194166; 0.09646265391194163; 3; 27.29177033333334; 81.87531100000001; 0.1612070459649248; 0.1612070459649248; 1202.874761361954; 10825.87285225759; 0.01555045948308101; 0.01555045948308101)
$rowsData[34] = @(3; 44.074633; 132.223899; 0.09646265391194166; 0.09646265391194163; 3; 40.05126633333333; 120.153799; 0.2365748448668895; 0.2365748448668895; 1765.244864826922; 15887.2037834423; 0.02282063738466605; 0.02282063738466604)
$rowsData[35] = @(3; 44.074633; 132.223899; 0.09646265391194166; 0.09646265391194163; 3; 40.10401433333334; 120.312043; 0.2368864167860689; 0.2368864167860689; 1767.569713568407; 15908.12742211566; 0.02285069243887453; 0.02285069243887452)
$rowsData[36] = @(3; 44.074633; 132.223899; 0.09646265391194166; 0.09646265391194163; 3; 29.373837; 88.121511; 0.1735053986454552; 0.1735053986454552; 1294.641085576821; 11651.76977019139; 0.01673679122139002; 0.01673679122139001)
$rowsData[37] = @(3; 44.074633; 132.223899; 0.09646265391194166; 0.09646265391194163; 3; 19.089213; 57.267639; 0.1127561752110562; 0.1127561752110562; 841.350057233829; 7572.15051510446; 0.01087675990581837; 0.01087675990581837)

foreach ($r in $rowsData.Keys) {
    $vals = $rowsData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value2 = $vals[$i]
    }
}

Write-Host "Applied NATMI App-Ncstn updates to rows 2-37"